# Deploy the implementation guide.
#
# - Metadata sheet: refresh the generated "Date" and "Contact" values, and
#   add a new "Jurisdiction" property row (pushing Description/Purpose/
#   Copyright/Immutable down by one row).
# - Rename the "Include from Ferlab.bio CodeS" sheet to "Include #0".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Refresh the generated metadata values.
$ws1.Range("B8").Value = "2024-10-02T15:04:17+00:00"
$ws1.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row for the "Jurisdiction" property right after "Contact",
# pushing the existing Description/Purpose/Copyright/Immutable rows down.
$ws1.Rows.Item(11).Insert()

# Match the look of the surrounding property rows (border/fill/alignment).
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# Rename the second sheet.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"
